$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# Insert 4 new rows after the current last data row (140). Inserting this way
# copies formatting down from the row above, so the new rows inherit the
# same per-column styles (date format in B, text format in C, the K column
# style, etc.) without having to fight with NumberFormat/Style plumbing.
$ws.Rows("141:144").Insert()

# Row 141 - a day with no recorded bet (mirrors the "empty" rows elsewhere,
# e.g. row 30/42/66/132: only A/B/C filled in, K kept blank but styled).
$ws.Range("A141").Value = 140
$ws.Range("B141").Value = 45246
$ws.Range("C141").Value = "2023-11-16"

# Row 142
$ws.Range("A142").Value = 141
$ws.Range("B142").Value = 45247
$ws.Range("C142").Value = "2023-11-17"
$ws.Range("D142").Value = 1
$ws.Range("E142").Value = 1.1299999999999999
$ws.Range("F142").Formula = '=H140'
$ws.Range("G142").Value = 55
$ws.Range("H142").Formula = '=F142+G142'
$ws.Range("I142").Value = "TENIS DE MESA"
$ws.Range("J142").Value = "SETKA CUP"
$ws.Range("K142").Formula = '=ROUND((H142/$F$31-1)*100, 3)+$K$29'

# Row 143
$ws.Range("A143").Value = 142
$ws.Range("B143").Value = 45247
$ws.Range("C143").Value = "2023-11-17"
$ws.Range("D143").Value = 1
$ws.Range("E143").Value = 1.155
$ws.Range("F143").Formula = '=H142'
$ws.Range("G143").Value = 40
$ws.Range("H143").Formula = '=F143+G143'
$ws.Range("I143").Value = "TENIS DE MESA"
$ws.Range("J143").Value = "LIGA PRO"
$ws.Range("K143").Formula = '=ROUND((H143/$F$31-1)*100, 3)+$K$29'

# Row 144
$ws.Range("A144").Value = 143
$ws.Range("B144").Value = 45247
$ws.Range("C144").Value = "2023-11-17"
$ws.Range("D144").Value = 1
$ws.Range("E144").Value = 1.07
$ws.Range("F144").Formula = '=H143'
$ws.Range("G144").Value = 21
$ws.Range("H144").Formula = '=F144+G144'
$ws.Range("I144").Value = "ESPORTS"
$ws.Range("J144").Value = "DOTA 2 ESL CHINA"
$ws.Range("K144").Formula = '=ROUND((H144/$F$31-1)*100, 3)+$K$29'

# Match the author's final selection (cell L142) recorded in the workbook view.
$ws.Range("L142").Select() | Out-Null
